$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-8, columns B (2) through J (10)
$data = @{
    2 = @(559, 511, 424, 391, 327, 280, 206, 92, 0)
    3 = @(71, 71, 65, 63, 60, 64, 76, 83, 89)
    4 = @(165, 136, 104, 95, 75, 60, 45, 24, 0)
    5 = @(74, 51, 38, 32, 29, 21, 17, 11, 0)
    6 = @(28, 59, 68, 68, 71, 73, 73, 76, 81)
    7 = @(52, 127, 259, 309, 397, 461, 542, 674, 790)
    8 = @(11, 5, 2, 2, 1, 1, 1, 0, 0)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # Column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
